$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly-entered 4mo_g (column G) belowground biomass readings for the rows
# that have a measurement. Rows without a matching value are left blank,
# same as in the source edit.
$values = @{
  5  = 1.516
  6  = 3.045
  7  = 1.562
  11 = 2.519
  12 = 1.016
  13 = 1.955
  14 = 1.204
  15 = 1.101
  16 = 1.333
  23 = 1.201
  24 = 1.517
  25 = 1.377
  29 = 0.695
  30 = 0.863
  31 = 0.931
  32 = 1.305
  33 = 0.505
  34 = 0.888
  38 = 1.265
  39 = 1.017
  40 = 1.53
  48 = 0.868
  49 = 0.953
  50 = 1.751
}

foreach ($row in $values.Keys) {
  $ws.Cells.Item($row, 7).Value = $values[$row]
}

# Match the author's final cursor position/scroll recorded in the sheet view.
$ws.Range("G51").Select()
